$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update final score (column K) values - each decreases by 0.3
$ws.Range("K2").Value = 54.6
$ws.Range("K3").Value = 50.6
$ws.Range("K4").Value = 46.4
$ws.Range("K5").Value = 45.2
$ws.Range("K6").Value = 36.6

# Update MACRO_SCORE (column N) values to the new constant for all rows
$ws.Range("N2:N6").Value = 50.60178744571824
